# Apply the numeric / structural updates described by the commit
# "calorimetry : scripts : tests : updated" to data.xlsx.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "equilibrium_concentrations" - tiny precision updates
# ---------------------------------------------------------------
$wsEq = $wb.Worksheets.Item("equilibrium_concentrations")
$wsEq.Range("A2").Value = 0.0000000000000000000000714307366483276
$wsEq.Range("F2").Value = 0.0000000000000000000999285700370759

$wsEq.Range("A4").Value = 0.0000000000000000000000714307360956512
$wsEq.Range("C4").Value = 0.00000584396962418032
$wsEq.Range("D4").Value = 0.00303747347236249
$wsEq.Range("F4").Value = 0.0000000000000000000999285692639044

$wsEq.Range("A6").Value = 0.0000000000000000000000714307360956517
$wsEq.Range("C6").Value = 0.00000584396962418031
$wsEq.Range("D6").Value = 0.00303747347236248
$wsEq.Range("F6").Value = 0.0000000000000000000999285692639058

# ---------------------------------------------------------------
# Sheet "heats_calculated" - tiny precision updates
# ---------------------------------------------------------------
$wsHeats = $wb.Worksheets.Item("heats_calculated")
$wsHeats.Range("I2").Value = -0.351639658307295
$wsHeats.Range("J2").Value = 0.0727445040872371

$wsHeats.Range("I3").Value = -0.231162839027781
$wsHeats.Range("J3").Value = 0.046451820397834

$wsHeats.Range("I4").Value = 0.60846327813482

# ---------------------------------------------------------------
# Sheet "enthalpies_calculated" - tiny precision update
# ---------------------------------------------------------------
$wsEnth = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnth.Range("C4").Value = 0.410638146694263

# ---------------------------------------------------------------
# Sheet "metrics" - precision updates + new RMSE row
# ---------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("metrics")
$wsMetrics.Range("B2").Value = -0.682948898173597
$wsMetrics.Range("B3").Value = 0.0831910329108941
$wsMetrics.Range("B4").Value = 0.0768819357757744

$wsMetrics.Range("A5").Value = "RMSE"
$wsMetrics.Range("B5").Value = 0.427127720274403
